$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 14 - "Fig 2: Main Menu Window" caption: remove italic formatting.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$capShape = $s14.Shapes.Item("Google Shape;135;p26")
$tf14 = $capShape.TextFrame
$savedText = $tf14.TextRange.Text
$tf14.DeleteText()
$tr14 = $tf14.TextRange
$tr14.Font.Italic = $false
$tr14.Text = $savedText

# ---------------------------------------------------------------------------
# Slides 17, 20, 23, 25, 27, 28 - reposition the figure-caption textboxes.
# EMU values below are divided by 12700 to get points (as PowerPoint's
# Shape.Left/Top expect); a tiny epsilon guards against float round-down
# when going back from points to EMU on save.
# ---------------------------------------------------------------------------
$emuPerPt = 12700
$eps = 0.5 / $emuPerPt

function Set-ShapeOffsetEmu {
    param($shape, [double]$xEmu, [double]$yEmu)
    $shape.Left = ($xEmu / $emuPerPt) + $eps
    $shape.Top  = ($yEmu / $emuPerPt) + $eps
}

$s17 = $p.Slides.Item(17)
Set-ShapeOffsetEmu ($s17.Shapes.Item("TextBox 4")) 2339787 4398352

$s20 = $p.Slides.Item(20)
Set-ShapeOffsetEmu ($s20.Shapes.Item("TextBox 2")) 1801906 4467642

$s23 = $p.Slides.Item(23)
Set-ShapeOffsetEmu ($s23.Shapes.Item("TextBox 2")) 2393575 4467642

$s25 = $p.Slides.Item(25)
Set-ShapeOffsetEmu ($s25.Shapes.Item("TextBox 1")) 1960736 4528306

$s27 = $p.Slides.Item(27)
Set-ShapeOffsetEmu ($s27.Shapes.Item("TextBox 1")) 2134224 4467642

$s28 = $p.Slides.Item(28)
Set-ShapeOffsetEmu ($s28.Shapes.Item("TextBox 1")) 2552736 4467642
